$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.346.04"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.698.86"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2.16"
$ws.Range("E5").Value = "  +12.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "236.70"
$ws.Range("E6").Value = "  -1.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "655.90"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.441"
$ws.Range("E8").Value = "  +3.13%  "

$ws.Range("E9").Value = "  +5.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.699.50"
$ws.Range("E11").Value = "  +0.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.97"
$ws.Range("E12").Value = "  -0.97%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000311"
$ws.Range("E13").Value = "  +15.08%  "

$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.83"
$ws.Range("E15").Value = "  -1.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.392.81"
$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.922.39"
$ws.Range("E17").Value = "  +0.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.91"
$ws.Range("E18").Value = "  -1.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.710.55"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.08"
$ws.Range("E20").Value = "  +1.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.89"
$ws.Range("E21").Value = "  -1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.539"
$ws.Range("E22").Value = "  +0.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "524.61"
$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("E24").Value = "  -2.00%  "

$ws.Range("E25").Value = "  +8.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "118.10"
$ws.Range("E26").Value = "  +15.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.93"
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.206"
$ws.Range("E28").Value = "  +22.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.42"
$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.74"
$ws.Range("E30").Value = "  +0.38%  "

$ws.Range("E31").Value = "  -1.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.188"
$ws.Range("E33").Value = "  +1.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.99"
$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.82"
$ws.Range("E35").Value = "  -3.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.596"
$ws.Range("E37").Value = "  -1.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "633.18"
$ws.Range("E38").Value = "  -3.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.73"
$ws.Range("E39").Value = "  -2.77%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("E41").Value = "  +14.18%  "

$ws.Range("E42").Value = "  +1.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.85"
$ws.Range("E43").Value = "  -4.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.24"
$ws.Range("E44").Value = "  +3.59%  "

$ws.Range("E45").Value = "  +0.21%  "

$ws.Range("E46").Value = "  -1.13%  "

$ws.Range("E47").Value = "  -1.26%  "

$ws.Range("E48").Value = "  +1.64%  "

$ws.Range("E49").Value = "  +1.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.64"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.32"
$ws.Range("E51").Value = "  +2.72%  "

